$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 4822.1577
$ws.Cells.Item(70, 9).Value = 4091.9
$ws.Cells.Item(70, 10).Value = 5633.5557
$ws.Cells.Item(70, 11).Value = 12275.7
$ws.Cells.Item(70, 12).Value = 16900.6671
$ws.Cells.Item(70, 13).Value = -12005.7
$ws.Cells.Item(70, 14).Value = -17440.6671

$ws.Cells.Item(73, 8).Value = 4822.1577
$ws.Cells.Item(73, 9).Value = 4091.9
$ws.Cells.Item(73, 10).Value = 5633.5557
$ws.Cells.Item(73, 11).Value = 12275.7
$ws.Cells.Item(73, 12).Value = 16900.6671
$ws.Cells.Item(73, 13).Value = -11339.7
$ws.Cells.Item(73, 14).Value = -18772.6671

$ws.Cells.Item(100, 8).Value = 2140.5908
$ws.Cells.Item(100, 9).Value = 1205.5
$ws.Cells.Item(100, 10).Value = 2919.8333
$ws.Cells.Item(100, 11).Value = 1205.5
$ws.Cells.Item(100, 12).Value = 2919.8333
$ws.Cells.Item(100, 13).Value = -664.5
$ws.Cells.Item(100, 14).Value = -4001.8333

$ws.Cells.Item(112, 8).Value = 1697.2106
$ws.Cells.Item(112, 9).Value = 1010.2222
$ws.Cells.Item(112, 10).Value = 2315.5
$ws.Cells.Item(112, 11).Value = 3030.6666
$ws.Cells.Item(112, 12).Value = 6946.5
$ws.Cells.Item(112, 13).Value = -1922.6666
$ws.Cells.Item(112, 14).Value = -9162.5

$ws.Cells.Item(116, 8).Value = 9222.727999999999
$ws.Cells.Item(116, 9).Value = 7823.6665
$ws.Cells.Item(116, 10).Value = 10901.6
$ws.Cells.Item(116, 11).Value = 7823.6665
$ws.Cells.Item(116, 12).Value = 10901.6
$ws.Cells.Item(116, 13).Value = -4381.6665
$ws.Cells.Item(116, 14).Value = -17785.6

$ws.Cells.Item(132, 8).Value = 1307.3823
$ws.Cells.Item(132, 9).Value = 1201.5938
$ws.Cells.Item(132, 11).Value = 3604.7814
$ws.Cells.Item(132, 13).Value = -1074.7814

$ws.Cells.Item(135, 8).Value = 533
$ws.Cells.Item(135, 9).Value = 521.8570999999999
$ws.Cells.Item(135, 11).Value = 4696.7139
$ws.Cells.Item(135, 13).Value = -2161.7139

$ws.Cells.Item(137, 8).Value = 2599.875
$ws.Cells.Item(137, 10).Value = 4718.5713
$ws.Cells.Item(137, 12).Value = 14155.7139
$ws.Cells.Item(137, 14).Value = -19255.7139

$ws.Cells.Item(138, 8).Value = 2874.6553
$ws.Cells.Item(138, 9).Value = 1992.2858
$ws.Cells.Item(138, 10).Value = 3698.2
$ws.Cells.Item(138, 11).Value = 5976.857400000001
$ws.Cells.Item(138, 12).Value = 11094.6
$ws.Cells.Item(138, 13).Value = -836.8574000000008
$ws.Cells.Item(138, 14).Value = -21374.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(31, 8).Value = 2212.2
$ws.Cells.Item(31, 9).Value = 2212.2
$ws.Cells.Item(31, 11).Value = 2212.2
$ws.Cells.Item(31, 13).Value = -1918.2

$ws.Cells.Item(32, 8).Value = 6243.8
$ws.Cells.Item(32, 9).Value = 5162.364
$ws.Cells.Item(32, 11).Value = 5162.364
$ws.Cells.Item(32, 13).Value = -4875.364

$ws.Cells.Item(45, 8).Value = 58827250
$ws.Cells.Item(45, 9).Value = 83334980
$ws.Cells.Item(45, 11).Value = 83334980
$ws.Cells.Item(45, 13).Value = -83334603

$ws.Cells.Item(61, 8).Value = 7082.8335
$ws.Cells.Item(61, 9).Value = 7082.8335
$ws.Cells.Item(61, 11).Value = 7082.8335
$ws.Cells.Item(61, 13).Value = -6870.8335

$ws.Cells.Item(63, 8).Value = 2524.75
$ws.Cells.Item(63, 9).Value = 2000
$ws.Cells.Item(63, 11).Value = 2000
$ws.Cells.Item(63, 13).Value = -1314

$ws.Cells.Item(66, 8).Value = 2524.75
$ws.Cells.Item(66, 9).Value = 2000
$ws.Cells.Item(66, 11).Value = 10000
$ws.Cells.Item(66, 13).Value = -6568

$ws.Cells.Item(132, 8).Value = 5873.174
$ws.Cells.Item(132, 9).Value = 1860.4
$ws.Cells.Item(132, 10).Value = 8959.923000000001
$ws.Cells.Item(132, 11).Value = 5581.200000000001
$ws.Cells.Item(132, 12).Value = 26879.769
$ws.Cells.Item(132, 13).Value = -3051.200000000001
$ws.Cells.Item(132, 14).Value = -31939.769

$ws.Cells.Item(136, 8).Value = 7082.8335
$ws.Cells.Item(136, 9).Value = 7082.8335
$ws.Cells.Item(136, 11).Value = 21248.5005
$ws.Cells.Item(136, 13).Value = -18698.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(11, 8).Value = 2091.5454
$ws.Cells.Item(11, 9).Value = 1938.5
$ws.Cells.Item(11, 10).Value = 2499.6667
$ws.Cells.Item(11, 11).Value = 1938.5
$ws.Cells.Item(11, 12).Value = 2499.6667
$ws.Cells.Item(11, 13).Value = -1798.5
$ws.Cells.Item(11, 14).Value = -2779.6667

$ws.Cells.Item(75, 8).Value = 28624.5
$ws.Cells.Item(75, 9).Value = 7249
$ws.Cells.Item(75, 11).Value = 7249
$ws.Cells.Item(75, 13).Value = -6313

$ws.Cells.Item(78, 8).Value = 28624.5
$ws.Cells.Item(78, 9).Value = 7249
$ws.Cells.Item(78, 11).Value = 21747
$ws.Cells.Item(78, 13).Value = -17067

$ws.Cells.Item(134, 8).Value = 2795.7
$ws.Cells.Item(134, 9).Value = 1456.3462
$ws.Cells.Item(134, 11).Value = 4369.0386
$ws.Cells.Item(134, 13).Value = -1834.0386

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1846.4
$ws.Cells.Item(22, 10).Value = 2519.4
$ws.Cells.Item(22, 12).Value = 2519.4
$ws.Cells.Item(22, 14).Value = -3219.4

$ws.Cells.Item(31, 8).Value = 21334.242
$ws.Cells.Item(31, 9).Value = 3272.2
$ws.Cells.Item(31, 11).Value = 3272.2
$ws.Cells.Item(31, 13).Value = -2977.2

$ws.Cells.Item(34, 8).Value = 21334.242
$ws.Cells.Item(34, 9).Value = 3272.2
$ws.Cells.Item(34, 11).Value = 3272.2
$ws.Cells.Item(34, 13).Value = -3070.2

$ws.Cells.Item(50, 8).Value = 28425
$ws.Cells.Item(50, 10).Value = 28425
$ws.Cells.Item(50, 12).Value = 28425
$ws.Cells.Item(50, 14).Value = -29675

$ws.Cells.Item(55, 8).Value = 23498
$ws.Cells.Item(55, 10).Value = 25000
$ws.Cells.Item(55, 12).Value = 25000
$ws.Cells.Item(55, 14).Value = -25630

$ws.Cells.Item(62, 8).Value = 6119.067
$ws.Cells.Item(62, 10).Value = 10867.5
$ws.Cells.Item(62, 12).Value = 10867.5
$ws.Cells.Item(62, 14).Value = -12115.5

$ws.Cells.Item(65, 8).Value = 6119.067
$ws.Cells.Item(65, 10).Value = 10867.5
$ws.Cells.Item(65, 12).Value = 54337.5
$ws.Cells.Item(65, 14).Value = -60577.5

$ws.Cells.Item(99, 8).Value = 4200
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 10).Value = 4200
$ws.Cells.Item(99, 11).Value = 0
$ws.Cells.Item(99, 12).Value = 4200
$ws.Cells.Item(99, 14).Value = -7196
$ws.Cells.Item(99, 13).ClearContents()

$ws.Cells.Item(126, 8).Value = 4200
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 4200
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 12).Value = 12600
$ws.Cells.Item(126, 14).Value = -17540
$ws.Cells.Item(126, 13).ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(6, 8).Value = 160.2
$ws.Cells.Item(6, 9).Value = 175.5
$ws.Cells.Item(6, 10).Value = 99
$ws.Cells.Item(6, 11).Value = 526.5
$ws.Cells.Item(6, 12).Value = 297
$ws.Cells.Item(6, 13).Value = -413.5
$ws.Cells.Item(6, 14).Value = -523

$ws.Cells.Item(26, 8).Value = 907.2727
$ws.Cells.Item(26, 9).Value = 873.1
$ws.Cells.Item(26, 11).Value = 2619.3
$ws.Cells.Item(26, 13).Value = -2331.3

$ws.Cells.Item(34, 8).Value = 5583.769
$ws.Cells.Item(34, 9).Value = 8234.143
$ws.Cells.Item(34, 10).Value = 2491.6667
$ws.Cells.Item(34, 11).Value = 24702.429
$ws.Cells.Item(34, 12).Value = 7475.000100000001
$ws.Cells.Item(34, 13).Value = -24618.429
$ws.Cells.Item(34, 14).Value = -7643.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(55, 8).Value = 521962.66
$ws.Cells.Item(55, 9).Value = 32998
$ws.Cells.Item(55, 10).Value = 766445
$ws.Cells.Item(55, 11).Value = 32998
$ws.Cells.Item(55, 12).Value = 766445
$ws.Cells.Item(55, 13).Value = -32671
$ws.Cells.Item(55, 14).Value = -767099

$ws.Cells.Item(59, 8).Value = 55000
$ws.Cells.Item(59, 9).Value = 0
$ws.Cells.Item(59, 11).Value = 0
$ws.Cells.Item(59, 13).ClearContents()

$ws.Cells.Item(126, 8).Value = 5501
$ws.Cells.Item(126, 9).Value = 1699.3334
$ws.Cells.Item(126, 10).Value = 6926.625
$ws.Cells.Item(126, 11).Value = 5098.0002
$ws.Cells.Item(126, 12).Value = 20779.875
$ws.Cells.Item(126, 13).Value = -2628.0002
$ws.Cells.Item(126, 14).Value = -25719.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 4121.4287
$ws.Cells.Item(46, 9).Value = 2966.3333
$ws.Cells.Item(46, 10).Value = 4987.75
$ws.Cells.Item(46, 11).Value = 2966.3333
$ws.Cells.Item(46, 12).Value = 4987.75
$ws.Cells.Item(46, 13).Value = -2778.3333
$ws.Cells.Item(46, 14).Value = -5363.75

$ws.Cells.Item(53, 8).Value = 51000
$ws.Cells.Item(53, 10).Value = 51000
$ws.Cells.Item(53, 12).Value = 51000
$ws.Cells.Item(53, 14).Value = -52036

$ws.Cells.Item(55, 8).Value = 2000484
$ws.Cells.Item(55, 9).Value = 3125421.8
$ws.Cells.Item(55, 10).Value = 594.8889
$ws.Cells.Item(55, 11).Value = 3125421.8
$ws.Cells.Item(55, 12).Value = 594.8889
$ws.Cells.Item(55, 13).Value = -3125248.8
$ws.Cells.Item(55, 14).Value = -940.8889

$ws.Cells.Item(82, 8).Value = 6955.2354
$ws.Cells.Item(82, 9).Value = 1317.5
$ws.Cells.Item(82, 10).Value = 11966.556
$ws.Cells.Item(82, 11).Value = 1317.5
$ws.Cells.Item(82, 12).Value = 11966.556
$ws.Cells.Item(82, 13).Value = -956.5
$ws.Cells.Item(82, 14).Value = -12688.556

$ws.Cells.Item(85, 8).Value = 6955.2354
$ws.Cells.Item(85, 9).Value = 1317.5
$ws.Cells.Item(85, 10).Value = 11966.556
$ws.Cells.Item(85, 11).Value = 1317.5
$ws.Cells.Item(85, 12).Value = 11966.556
$ws.Cells.Item(85, 13).Value = -69.5
$ws.Cells.Item(85, 14).Value = -14462.556

$ws.Cells.Item(136, 8).Value = 4882.8
$ws.Cells.Item(136, 9).Value = 2220.8333
$ws.Cells.Item(136, 11).Value = 6662.499899999999
$ws.Cells.Item(136, 13).Value = -4112.499899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(45, 8).Value = 14590.5
$ws.Cells.Item(45, 9).Value = 0
$ws.Cells.Item(45, 10).Value = 14590.5
$ws.Cells.Item(45, 11).Value = 0
$ws.Cells.Item(45, 12).Value = 14590.5
$ws.Cells.Item(45, 14).Value = -15572.5
$ws.Cells.Item(45, 13).ClearContents()

$ws.Cells.Item(132, 8).Value = 3350.7273
$ws.Cells.Item(132, 9).Value = 3359.639
$ws.Cells.Item(132, 11).Value = 10078.917
$ws.Cells.Item(132, 13).Value = -7548.917000000001
